$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(49)
Write-Output ("Text: " + $p.Range.Text)
Write-Output ("RightIndent: " + $p.Range.ParagraphFormat.RightIndent)
Write-Output ("LeftIndent: " + $p.Range.ParagraphFormat.LeftIndent)
Write-Output ("FirstLineIndent: " + $p.Range.ParagraphFormat.FirstLineIndent)
